$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.493.04'
$ws.Range("E2").Value = '  +0.93%  '
$ws.Range("D3").Value = '2.985.83'
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '381.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.82%  '
$ws.Range("E7").Value = '  +1.18%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +0.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.68'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("E11").Value = '  -0.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0859'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.08%  '
$ws.Range("D13").Value = '3.451.59'
$ws.Range("E13").Value = '  +1.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.80'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.88%  '
$ws.Range("D16").Value = '2.980.64'
$ws.Range("E16").Value = '  +1.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '11.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.00%  '
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("D19").Value = '51.481.34'
$ws.Range("E19").Value = '  +0.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("D22").Value = '0.0₃0964'
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '266.88'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.67%  '
$ws.Range("E25").Value = '  +2.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.73%  '
$ws.Range("E27").Value = '  -1.50%  '
$ws.Range("E28").Value = '  +3.12%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '26.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.02%  '
$ws.Range("E31").Value = '  -0.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.41'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.83'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.44%  '
$ws.Range("E34").Value = '  +1.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.74%  '
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.29'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.85'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.58'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.52%  '
$ws.Range("E41").Value = '  +0.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.85'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '127.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.93%  '
$ws.Range("E44").Value = '  +12.81%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.43'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("E47").Value = '  +0.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.36'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.71%  '
$ws.Range("D49").Value = '2.028.66'
$ws.Range("E49").Value = '  +2.15%  '
$ws.Range("D50").Value = '3.280.67'
$ws.Range("E50").Value = '  +1.26%  '
$ws.Range("E51").Value = '  +2.01%  '
